$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New curriculum rows (semester-selector feature): replace rows 2-4 and
# append through row 14. Columns: A name, B teacher, C room, D course code,
# E section no., F weeks, G weekday (number), H credits.
# D/E/H are numeric-looking text (e.g. "05", "01") that Excel would silently
# coerce to numbers (dropping leading zeros) on plain assignment, so those
# three columns are written with a temporary text format.
$courseData = @(
    @("体育-3", "ty9", "江安体育场体育场2号", "888006010", "09", "2-17周", 3, "1"),
    @("形势与政策-3", "刘辉", "江安综合楼C座C306", "107117000", "24", "3-9周单", 5, "0"),
    @("马克思主义基本原理概论", "陈伟", "江安综合楼C座C403", "107021030", "36", "2-18周", 1, "3"),
    @("JAVA程序设计", "李琳", "江安一教B座B203", "304007030", "05", "2-17周", 5, "3"),
    @("大学英语（创意阅读）-3", "徐光源", "江安一教A座A203", "105395020", "230", "2-18周", 3, "2"),
    @("物联网工程导论", "桑永胜", "江安一教A座A506", "304085010", "01", "2-9周", 1, "1"),
    @("物联网传感器原理及应用", "时宏伟", "江安一教C座C504", "304211030", "01", "2-9周", 2, "3"),
    @("计算机组成原理实验", "蒋欣荣", "江安实验室二基楼B310B", "304037010", "06", "11-15周", 2, "1"),
    @("离散数学", "陈瑜", "江安综合楼C座C303", "304156050", "07", "2-18周", 1, "5"),
    @("离散数学", "陈瑜", "江安一教B座B201", "304156050", "07", "2-18周", 4, "5"),
    @("计算机组成原理", "蒋欣荣", "江安综合楼C座C408", "304036030", "06", "2-17周", 2, "3"),
    @("数据结构与算法分析课程设计", "周欣", "江安实验室二基楼B304", "304046010", "07", "6-15周", 4, "1"),
    @("数据结构与算法分析", "周欣", "江安一教C座C406", "304045030", "07", "2-17周", 4, "3")
)

$startRow = 2
for ($i = 0; $i -lt $courseData.Count; $i++) {
    $row = $startRow + $i
    $rec = $courseData[$i]

    $ws.Cells.Item($row, 1).Value = $rec[0]   # A: course name
    $ws.Cells.Item($row, 2).Value = $rec[1]   # B: teacher
    $ws.Cells.Item($row, 3).Value = $rec[2]   # C: room

    # D, E, H: force text so values like "05"/"01" keep their leading zero
    # instead of being auto-converted to numbers.
    foreach ($col in @(4, 5, 8)) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.NumberFormat = "@"
        $cell.Value = $rec[$col - 1]
        $cell.NumberFormat = "General"
    }

    $ws.Cells.Item($row, 6).Value = $rec[5]   # F: weeks
    $ws.Cells.Item($row, 7).Value = $rec[6]   # G: weekday (numeric)
}

